$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B2").Value = "C01号直流"
$ws.Range("C2").Value = "2025-01-25 13:46:36"
$ws.Range("D2").Value = 45948.29740740741

$ws.Range("A3").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B3").Value = "C02号直流"
$ws.Range("C3").Value = "2025-01-25 17:13:47"
$ws.Range("D3").Value = 45948.29740740741

$ws.Range("A4").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B4").Value = "C03号直流"
$ws.Range("C4").Value = "2025-01-25 14:14:24"
$ws.Range("D4").Value = 45948.29740740741

$ws.Range("A5").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B5").Value = "C04号直流"
$ws.Range("C5").Value = "2025-01-25 06:24:40"
$ws.Range("D5").Value = 45948.29740740741

$ws.Range("A6").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B6").Value = "C05号直流"
$ws.Range("C6").Value = "2025-01-25 16:01:40"
$ws.Range("D6").Value = 45948.29740740741

$ws.Range("A7").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B7").Value = "D01号直流"
$ws.Range("C7").Value = "2025-01-25 18:30:24"
$ws.Range("D7").Value = 45948.29740740741

$ws.Range("A8").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B8").Value = "D02号直流"
$ws.Range("C8").Value = "2025-01-25 15:39:19"
$ws.Range("D8").Value = 45948.29740740741

$ws.Range("A9").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B9").Value = "D03号直流"
$ws.Range("C9").Value = "2025-01-25 16:09:35"
$ws.Range("D9").Value = 45948.29740740741

$ws.Range("A10").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B10").Value = "D04号直流"
$ws.Range("C10").Value = "2025-01-25 18:29:02"
$ws.Range("D10").Value = 45948.29740740741

$ws.Range("A11").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B11").Value = "D05号直流"
$ws.Range("C11").Value = "2025-01-25 18:27:29"
$ws.Range("D11").Value = 45948.29740740741

$ws.Range("A12").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B12").Value = "E01号直流"
$ws.Range("C12").Value = "2025-01-25 15:22:58"
$ws.Range("D12").Value = 45948.29740740741

$ws.Range("A13").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B13").Value = "E02号直流"
$ws.Range("C13").Value = "2025-01-25 16:45:57"
$ws.Range("D13").Value = 45948.29740740741

$ws.Range("A14").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B14").Value = "E03号直流"
$ws.Range("C14").Value = "2025-01-25 02:54:59"
$ws.Range("D14").Value = 45948.29740740741

$ws.Range("A15").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B15").Value = "E04号直流"
$ws.Range("C15").Value = "2025-01-25 17:08:37"
$ws.Range("D15").Value = 45948.29740740741

$ws.Range("A16").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B16").Value = "004B号直流"
$ws.Range("C16").Value = "2025-02-19 00:26:27"
$ws.Range("D16").Value = 45948.29740740741

$ws.Range("A17").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B17").Value = "701号直流"
$ws.Range("C17").Value = 45927.457337962966
$ws.Range("D17").Value = 45948.29740740741

$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "101号直流"
$ws.Range("C18").Value = 45935.0465625
$ws.Range("D18").Value = 45948.29740740741

$ws.Range("A19").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B19").Value = "406号直流"
$ws.Range("C19").Value = 45943.02091435185
$ws.Range("D19").Value = 45948.29740740741

$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "604号直流"
$ws.Range("C20").Value = 45946.1093287037
$ws.Range("D20").Value = 45948.29740740741

$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "505号直流"
$ws.Range("C21").Value = 45946.557071759256
$ws.Range("D21").Value = 45948.29740740741

$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "702号直流"
$ws.Range("C22").Value = 45947.04481481481
$ws.Range("D22").Value = 45948.29740740741

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "801号直流"
$ws.Range("C23").Value = 45947.17630787037
$ws.Range("D23").Value = 45948.29740740741

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "705号直流"
$ws.Range("C24").Value = 45947.264085648145
$ws.Range("D24").Value = 45948.29740740741

$ws.Range("A25").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B25").Value = "108号直流"
$ws.Range("C25").Value = 45947.51122685185
$ws.Range("D25").Value = 45948.29740740741

$ws.Range("A26").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B26").Value = "210号直流"
$ws.Range("C26").Value = 45947.52209490741
$ws.Range("D26").Value = 45948.29740740741

$ws.Range("A27").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B27").Value = "112号直流"
$ws.Range("C27").Value = 45947.52819444444
$ws.Range("D27").Value = 45948.29740740741

$ws.Range("A28").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B28").Value = "107号直流"
$ws.Range("C28").Value = 45947.550729166665
$ws.Range("D28").Value = 45948.29740740741

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "404号直流"
$ws.Range("C29").Value = 45947.55625
$ws.Range("D29").Value = 45948.29740740741

$ws.Range("A30").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B30").Value = "110号直流"
$ws.Range("C30").Value = 45947.5584375
$ws.Range("D30").Value = 45948.29740740741

$ws.Range("A31").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B31").Value = "201号直流"
$ws.Range("C31").Value = 45947.572592592594
$ws.Range("D31").Value = 45948.29740740741

$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "703号直流"
$ws.Range("C32").Value = 45947.581458333334
$ws.Range("D32").Value = 45948.29740740741

$ws.Range("A33").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B33").Value = "103号直流"
$ws.Range("C33").Value = 45947.581712962965
$ws.Range("D33").Value = 45948.29740740741

$ws.Range("A34").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B34").Value = "106号直流"
$ws.Range("C34").Value = 45947.58215277778
$ws.Range("D34").Value = 45948.29740740741

$ws.Range("A35").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value = "002A号直流"
$ws.Range("C35").Value = 45947.58777777778
$ws.Range("D35").Value = 45948.29740740741

$ws.Range("A36").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B36").Value = "504号直流"
$ws.Range("C36").Value = 45947.59570601852
$ws.Range("D36").Value = 45948.29740740741

$ws.Range("A37").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B37").Value = "105号直流"
$ws.Range("C37").Value = 45947.598703703705
$ws.Range("D37").Value = 45948.29740740741

$ws.Range("A38").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B38").Value = "804号直流"
$ws.Range("C38").Value = 45947.62819444444
$ws.Range("D38").Value = 45948.29740740741

$ws.Range("A39").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B39").Value = "307号直流"
$ws.Range("C39").Value = 45947.77392361111
$ws.Range("D39").Value = 45948.29740740741

$ws.Range("A40:E42").ClearContents()

$ws.Range("H21").Select()